$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values (prices) remain plain text, matching the
# workbook's original inlineStr storage, instead of being coerced to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.66"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "25.17"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.046"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05596"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.555"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.012"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8153"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8350"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1336"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06964"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03250"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02836"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09393"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001510"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0005969"
$ws.Range("E16").Value = "15OneONE"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006074"
$ws.Range("E17").Value = "16TigerCashTCH"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.499"
$ws.Range("E18").Value = "17LEOLEO"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.092"
$ws.Range("E19").Value = "18BTSETokenBTSE"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3188"
$ws.Range("E20").Value = "19BitpandaEcosystemTokenBEST"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.739"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04690"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001244"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009698"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001940"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03670"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006199"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1050"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002600"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008218"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1800"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002016"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
